# Generate Report for Handback
# Update the timestamp text values that record when the HO/Handback XLIFF
# files were (re)generated for the "03be178d-..." row across the
# Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" for 03be178d row ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Replace("2016-08-18 10:48:22", "2016-08-18 10:49:09")

# --- zh-cn sheet: Correspond Handoff / Handback DateTime for 03be178d row ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Cells.Replace("2016-08-18 10:48:17", "2016-08-18 10:48:59")
$wsZhCn.Cells.Replace("2016-08-18 10:48:35", "2016-08-18 10:49:29")

# --- de-de sheet: Correspond Handback DateTime for 03be178d row ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Cells.Replace("2016-08-18 10:48:42", "2016-08-18 10:49:36")
